{"js": "// The edit collapses the document's first three paragraphs:\n//   1. \"Subtitle: One-Stop Date Suggestions\" (+ the _GoBack bookmark)\n//   2. an empty paragraph\n//   3. \"Description: NoWut is an application...\" split across several\n//      runs (with spell-check <w:proofErr/> wrappers around \"NoWut\")\n// into a single paragraph that keeps the _GoBack bookmark followed by\n// ONE run holding the full \"Description: ...\" text (runs merged, the\n// \"Subtitle\" text removed, proofErr markers gone). The final trailing\n// empty paragraph is left untouched.\n\nconst FULL_TEXT =\n  \"Description: NoWut is an application that, in its most simple form, \" +\n  \"provides users with date suggestions.  It is designed to be completely \" +\n  \"compatible with user preferences.  Users can input various constraints \" +\n  \"in the form of filters in order to find activities that are most \" +\n  \"suited to their tastes.  The application takes the filters set by the \" +\n  \"user and searches through the google API for any locations that \" +\n  \"match.  Results are then provided to the user, both in map view and \" +\n  \"list view.  Users may filter the results by distance, cost and other \" +\n  \"various finely tuned options to be as accurate as they desire.  Users \" +\n  \"can also make accounts in order to store their search results in case \" +\n  \"they would like to repeat a previous search.  This application can \" +\n  \"find all sorts of entertainment, food, activities and much more.  \" +\n  \"With NoWut, a user will never have to worry \\u201Cwhat\\u2019s next?\\u201D.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// paragraphs.items[0] => \"Subtitle: ...\" (holds the _GoBack bookmark)\n// paragraphs.items[1] => empty paragraph\n// paragraphs.items[2] => \"Description: ...\" (multiple runs)\n// paragraphs.items[3] => trailing empty paragraph (untouched)\nconst subtitlePara = paragraphs.items[0];\n\n// Rebuild paragraph 0's content in one shot: bookmark tags first, then a\n// single run with the full merged description text - this matches the\n// element order Word itself keeps bookmarks that sit at a paragraph's\n// end anchored ahead of freshly-inserted content.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t>' + FULL_TEXT + '</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nsubtitlePara.getRange(\"Whole\").insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-query: paragraph 0 now carries the merged text; paragraph 1 is\n// still the old empty spacer paragraph; paragraph 2 is still the old\n// multi-run \"Description: ...\" paragraph (now redundant) and must go.\nlet current = context.document.body.paragraphs;\ncurrent.load(\"items\");\nawait context.sync();\n\ncurrent.items[2].delete(); // old \"Description: ...\" paragraph\nawait context.sync();\n\ncurrent = context.document.body.paragraphs;\ncurrent.load(\"items\");\nawait context.sync();\n\ncurrent.items[1].delete(); // old empty spacer paragraph\nawait context.sync();\n", "ps1": "# The edit collapses the document's first three paragraphs:\n#   1. \"Subtitle: One-Stop Date Suggestions\" (+ the _GoBack bookmark)\n#   2. an empty paragraph\n#   3. \"Description: NoWut is an application...\" split across several\n#      runs (with spell-check proofErr wrappers around \"NoWut\")\n# into a single paragraph that keeps the _GoBack bookmark followed by\n# ONE run holding the full \"Description: ...\" text (runs merged, the\n# \"Subtitle\" text removed, proofErr markers gone). The final trailing\n# empty paragraph is left untouched.\n\n$d = $word.ActiveDocument\n\n$fullText = \"Description: NoWut is an application that, in its most simple form, provides users with date suggestions.  It is designed to be completely compatible with user preferences.  Users can input various constraints in the form of filters in order to find activities that are most suited to their tastes.  The application takes the filters set by the user and searches through the google API for any locations that match.  Results are then provided to the user, both in map view and list view.  Users may filter the results by distance, cost and other various finely tuned options to be as accurate as they desire.  Users can also make accounts in order to store their search results in case they would like to repeat a previous search.  This application can find all sorts of entertainment, food, activities and much more.  With NoWut, a user will never have to worry \" + [char]0x201C + \"what\" + [char]0x2019 + \"s next?\" + [char]0x201D + \".\"\n\n# Paragraph 1 holds the _GoBack bookmark (at the end of its range) plus\n# the \"Subtitle: ...\" run. Rebuild it in one shot via InsertXML so the\n# bookmark tags land ahead of the new run, matching the target markup.\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>' + $fullText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$d.Paragraphs.Item(1).Range.InsertXML($ooxml)\n\n# Paragraph 3 is now the old, now-redundant \"Description: ...\" paragraph\n# (made of several runs); paragraph 2 is still the old empty spacer\n# paragraph. Remove both, leaving the merged paragraph followed by the\n# untouched trailing empty paragraph.\n$d.Paragraphs.Item(3).Range.Delete()\n$d.Paragraphs.Item(2).Range.Delete()\n"}
